# ---------------------------------------------------------------------------
# ppgcs_prodpubl_2017_2024.xlsx - "atualizacao dados docentes e discentes
# Fiocruz Ceara"
#
# The Power BI export backing this sheet was refreshed:
#   - Two advisors gained extra 2024 (column I) points:
#       Carla Freire Celedonio Fernandes : I 170 -> 230 (J 1465 -> 1525)
#       Anna Carolina Machado Marinho    : I  90 -> 150 (J  270 ->  330)
#   - "Antonio Marcos Aires Barbosa" dropped out of both the Autor table
#     (columns A:J) and the colab_disc reference list (column K), so the
#     sheet shrank from 40 to 39 rows (dimension A1:K40 -> A1:K39).
#   - Because two totals changed, the table (sorted descending by
#     "Soma de Pontos", column J) was re-sorted.
#
# Row 1 (headers) is untouched. Rows 2-39 below hold the resulting table,
# row 40 is removed.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,('Jaime Ribeiro Filho', 575, 320, 200, 1060, 955, 865, 1155, 200, 5330, 'Alice Paula Di Sabatino Guimarães')
  ,('Márcio Flávio Moura de Araújo', 285, 570, 315, 385, 475, 260, 380, 80, 2750, 'Ana Cláudia de Araújo Teixeira')
  ,('Fabio Miyajima', 350, 270, 330, 270, 530, 770, 80, 0, 2600, 'Ana Camila Oliveira Alves')
  ,('Anya Pimentel Gomes Fernandes Vieira Meyer', 100, 310, 230, 355, 755, 300, 445, 100, 2595, 'Adriana Costa Bacelo')
  ,('Roberto Wagner Júnior Freire de Freitas', 175, 390, 320, 365, 580, 80, 190, 20, 2120, 'Anna Carolina Machado Marinho')
  ,('Roberto Nicolete', 160, 80, 55, 270, 335, 510, 285, 80, 1775, 'Anya Pimentel Gomes Fernandes Vieira Meyer')
  ,('Ivana Cristina de Holanda Cunha Barreto', 335, 180, 65, 150, 370, 340, 200, 0, 1640, 'Carla Freire Celedonio Fernandes')
  ,('Carla Freire Celedonio Fernandes', 410, 170, 100, 220, 320, 60, 15, 230, 1525, 'Claudia Stutz Zubieta')
  ,('Caroline Pereira Bittencourt Passaes', 90, 90, 110, 360, 490, 170, 90, 90, 1490, 'Clarissa Romero Teixeira')
  ,('João Hermínio Martins da Silva', 310, 0, 210, 120, 55, 520, 200, 0, 1415, 'Dayane Alves Costa')
  ,('Luiz Odorico Monteiro de Andrade', 170, 130, 15, 150, 420, 260, 140, 40, 1325, 'Donat Alexander de Chapeaurouge')
  ,('José Luís Passos Cordeiro', 0, 400, 100, 375, 90, 130, 80, 15, 1190, 'Eduardo Ruback dos Santos')
  ,('Sharmênia de Araújo Soares Nuto', 80, 55, 50, 240, 385, 115, 170, 60, 1155, 'Fabio Miyajima')
  ,('Fernando Ferreira Carneiro', 170, 205, 40, 130, 130, 100, 180, 180, 1135, 'Fernando Braga Stehling Dias')
  ,('Marcos Roberto Lourenzoni', 110, 90, 180, 190, 120, 210, 100, 0, 1000, 'Fernando Ferreira Carneiro')
  ,('Fernando Braga Stehling Dias', 90, 90, 90, 80, 260, 270, 80, 0, 960, 'Galba Freire Moita')
  ,('Gilvan Pessoa Furtado', 100, 180, 110, 0, 270, 100, 160, 0, 920, 'Giovanny Augusto Camacho Antevere Mazzarotto')
  ,('Clarissa Romero Teixeira', 0, 320, 0, 80, 0, 80, 180, 80, 740, 'Gilvan Pessoa Furtado')
  ,('Maximiliano Ponte', 270, 20, 75, 130, 90, 0, 95, 0, 680, 'Ivana Cristina de Holanda Cunha Barreto')
  ,('Vanira Matos Pessoa', 150, 80, 40, 20, 125, 110, 150, 0, 675, 'Jaime Ribeiro Filho')
  ,('Marcela Helena Gambim Fonseca', 0, 90, 20, 0, 170, 170, 180, 10, 640, 'João Hermínio Martins da Silva')
  ,('Regis Bernardo Brandim Gomes', 0, 0, 0, 40, 20, 160, 180, 200, 600, 'José Luís Passos Cordeiro')
  ,('Eduardo Ruback dos Santos', 0, 0, 0, 0, 0, 270, 0, 90, 360, 'Luiz Odorico Monteiro de Andrade')
  ,('Anna Carolina Machado Marinho', 0, 15, 0, 0, 150, 0, 15, 150, 330, 'Marcela Helena Gambim Fonseca')
  ,('Marlos de Medeiros Chaves', 0, 90, 0, 40, 80, 90, 0, 0, 300, 'Marcos Roberto Lourenzoni')
  ,('Adriana Costa Bacelo', 50, 0, 110, 0, 10, 35, 0, 80, 285, 'Márcio Flávio Moura de Araújo')
  ,('Venúcia Bruna Magalhães Pereira', 80, 80, 60, 0, 40, 0, 0, 0, 260, 'Margareth Borges Coutinho Gallo')
  ,('Claudia Stutz Zubieta', 0, 0, 0, 0, 0, 90, 90, 80, 260, 'Marlos de Medeiros Chaves')
  ,('Galba Freire Moita', 0, 20, 95, 40, 25, 5, 0, 0, 185, 'Maximiliano Ponte')
  ,('Alice Paula Di Sabatino Guimarães', 0, 0, 0, 0, 170, 0, 0, 0, 170, 'Raphael Trevizani')
  ,('Dayane Alves Costa', 90, 0, 0, 0, 40, 0, 20, 0, 150, 'Regis Bernardo Brandim Gomes')
  ,('Ana Camila Oliveira Alves', 0, 0, 0, 0, 60, 60, 0, 0, 120, 'Roberto Nicolete')
  ,('Ana Cláudia de Araújo Teixeira', 0, 0, 0, 0, 0, 20, 90, 0, 110, 'Roberto Wagner Júnior Freire de Freitas')
  ,('Margareth Borges Coutinho Gallo', 0, 0, 0, 0, 100, 0, 0, 0, 100, 'Sharmênia de Araújo Soares Nuto')
  ,('Giovanny Augusto Camacho Antevere Mazzarotto', 80, 0, 0, 0, 5, 0, 15, 0, 100, 'Vanira Matos Pessoa')
  ,('Raphael Trevizani', 90, 0, 0, 0, 0, 0, 0, 0, 90, 'Venúcia Bruna Magalhães Pereira')
  ,('Donat Alexander de Chapeaurouge', 0, 80, 0, 0, 0, 0, 0, 0, 80, 'Fernanda Savicki de Almeida')
  ,('Fernanda Savicki de Almeida', 0, 0, 0, 0, 0, 0, 5, 0, 5, 'Caroline Pereira Bittencourt Passaes')
)

# Row 1 headers are unchanged - only rewrite rows 2..39 (11 columns: A..K).
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    for ($c = 2; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $ws.Cells.Item($r, 11).Value = $row[10]
    $r++
}

# The sheet shrank from 40 to 39 rows (Antonio Marcos Aires Barbosa removed);
# drop the now-stale last row and let Excel recompute the used range.
$ws.Rows.Item(40).Delete()
